$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'FL8222'
$ws.Cells.Item(2, 2).Value = 'website'
$ws.Cells.Item(2, 3).Value = 'premium_economy'
$ws.Cells.Item(2, 4).Value = 9
$ws.Cells.Item(2, 5).Value = 1684.31
$ws.Cells.Item(2, 6).Value = 412.16
$ws.Cells.Item(2, 7).Value = 750
$ws.Cells.Item(2, 8).Value = 2846.47
$ws.Cells.Item(2, 9).Value = 316.27
$ws.Cells.Item(2, 10).Value = 142.32

# Row 3
$ws.Cells.Item(3, 1).Value = 'FL1013'
$ws.Cells.Item(3, 2).Value = 'travel_agent'
$ws.Cells.Item(3, 3).Value = 'premium_economy'
$ws.Cells.Item(3, 4).Value = 13
$ws.Cells.Item(3, 5).Value = 7398.05
$ws.Cells.Item(3, 6).Value = 784.88
$ws.Cells.Item(3, 7).Value = 1100
$ws.Cells.Item(3, 8).Value = 9282.93
$ws.Cells.Item(3, 9).Value = 714.0700000000001
$ws.Cells.Item(3, 10).Value = 464.15

# Row 4
$ws.Cells.Item(4, 1).Value = 'FL8121'
$ws.Cells.Item(4, 2).Value = 'app'
$ws.Cells.Item(4, 3).Value = 'premium_economy'
$ws.Cells.Item(4, 4).Value = 7
$ws.Cells.Item(4, 5).Value = 1077.11
$ws.Cells.Item(4, 6).Value = 420.55
$ws.Cells.Item(4, 7).Value = 700
$ws.Cells.Item(4, 8).Value = 2197.66
$ws.Cells.Item(4, 9).Value = 313.95
$ws.Cells.Item(4, 10).Value = 109.88

# Row 5
$ws.Cells.Item(5, 1).Value = 'FL1541'
$ws.Cells.Item(5, 2).Value = 'app'
$ws.Cells.Item(5, 3).Value = 'premium_economy'
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 6876.5
$ws.Cells.Item(5, 6).Value = 399.17
$ws.Cells.Item(5, 7).Value = 350
$ws.Cells.Item(5, 8).Value = 7625.67
$ws.Cells.Item(5, 9).Value = 1270.95
$ws.Cells.Item(5, 10).Value = 381.28

# Row 6
$ws.Cells.Item(6, 1).Value = 'FL8464'
$ws.Cells.Item(6, 2).Value = 'travel_agent'
$ws.Cells.Item(6, 3).Value = 'economy'
$ws.Cells.Item(6, 4).Value = 26
$ws.Cells.Item(6, 5).Value = 5370.77
$ws.Cells.Item(6, 6).Value = 1100.34
$ws.Cells.Item(6, 7).Value = 2250
$ws.Cells.Item(6, 8).Value = 8721.110000000001
$ws.Cells.Item(6, 9).Value = 335.43
$ws.Cells.Item(6, 10).Value = 436.06

# Row 7
$ws.Cells.Item(7, 1).Value = 'FL7466'
$ws.Cells.Item(7, 2).Value = 'website'
$ws.Cells.Item(7, 3).Value = 'economy'
$ws.Cells.Item(7, 4).Value = 31
$ws.Cells.Item(7, 5).Value = 6301.85
$ws.Cells.Item(7, 6).Value = 1072.56
$ws.Cells.Item(7, 7).Value = 2350
$ws.Cells.Item(7, 8).Value = 9724.41
$ws.Cells.Item(7, 9).Value = 313.69
$ws.Cells.Item(7, 10).Value = 486.22

# Row 8
$ws.Cells.Item(8, 1).Value = 'FL2900'
$ws.Cells.Item(8, 2).Value = 'website'
$ws.Cells.Item(8, 3).Value = 'premium_economy'
$ws.Cells.Item(8, 4).Value = 8
$ws.Cells.Item(8, 5).Value = 1379.34
$ws.Cells.Item(8, 6).Value = 346.88
$ws.Cells.Item(8, 7).Value = 350
$ws.Cells.Item(8, 8).Value = 2076.22
$ws.Cells.Item(8, 9).Value = 259.53
$ws.Cells.Item(8, 10).Value = 103.81

# Row 9
$ws.Cells.Item(9, 1).Value = 'FL1740'
$ws.Cells.Item(9, 2).Value = 'app'
$ws.Cells.Item(9, 3).Value = 'premium_economy'
$ws.Cells.Item(9, 4).Value = 6
$ws.Cells.Item(9, 5).Value = 1215.49
$ws.Cells.Item(9, 6).Value = 396.16
$ws.Cells.Item(9, 7).Value = 600
$ws.Cells.Item(9, 8).Value = 2211.65
$ws.Cells.Item(9, 9).Value = 368.61
$ws.Cells.Item(9, 10).Value = 110.58

# Row 10
$ws.Cells.Item(10, 1).Value = 'FL2948'
$ws.Cells.Item(10, 2).Value = 'website'
$ws.Cells.Item(10, 3).Value = 'business'
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 1414.76
$ws.Cells.Item(10, 6).Value = 267.26
$ws.Cells.Item(10, 7).Value = 150
$ws.Cells.Item(10, 8).Value = 1832.02
$ws.Cells.Item(10, 9).Value = 366.4
$ws.Cells.Item(10, 10).Value = 91.59999999999999

# Row 11
$ws.Cells.Item(11, 1).Value = 'FL5986'
$ws.Cells.Item(11, 2).Value = 'app'
$ws.Cells.Item(11, 3).Value = 'business'
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 1158.83
$ws.Cells.Item(11, 6).Value = 147.64
$ws.Cells.Item(11, 7).Value = 500
$ws.Cells.Item(11, 8).Value = 1806.47
$ws.Cells.Item(11, 9).Value = 361.29
$ws.Cells.Item(11, 10).Value = 90.31999999999999

# Row 12
$ws.Cells.Item(12, 1).Value = 'FL3066'
$ws.Cells.Item(12, 2).Value = 'travel_agent'
$ws.Cells.Item(12, 3).Value = 'business'
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = 1263.64
$ws.Cells.Item(12, 6).Value = 184.66
$ws.Cells.Item(12, 7).Value = 400
$ws.Cells.Item(12, 8).Value = 1848.3
$ws.Cells.Item(12, 9).Value = 369.66
$ws.Cells.Item(12, 10).Value = 92.42

# Row 13
$ws.Cells.Item(13, 1).Value = 'FL2625'
$ws.Cells.Item(13, 2).Value = 'website'
$ws.Cells.Item(13, 3).Value = 'economy'
$ws.Cells.Item(13, 4).Value = 28
$ws.Cells.Item(13, 5).Value = 3499.14
$ws.Cells.Item(13, 6).Value = 1302.44
$ws.Cells.Item(13, 7).Value = 2050
$ws.Cells.Item(13, 8).Value = 6851.58
$ws.Cells.Item(13, 9).Value = 244.7
$ws.Cells.Item(13, 10).Value = 342.58

# Row 14
$ws.Cells.Item(14, 1).Value = 'FL4736'
$ws.Cells.Item(14, 2).Value = 'app'
$ws.Cells.Item(14, 3).Value = 'economy'
$ws.Cells.Item(14, 4).Value = 45
$ws.Cells.Item(14, 5).Value = 5610.04
$ws.Cells.Item(14, 6).Value = 2143.88
$ws.Cells.Item(14, 7).Value = 3550
$ws.Cells.Item(14, 8).Value = 11303.92
$ws.Cells.Item(14, 9).Value = 251.2
$ws.Cells.Item(14, 10).Value = 565.2

# Row 15
$ws.Cells.Item(15, 1).Value = 'FL9612'
$ws.Cells.Item(15, 2).Value = 'website'
$ws.Cells.Item(15, 3).Value = 'premium_economy'
$ws.Cells.Item(15, 4).Value = 14
$ws.Cells.Item(15, 5).Value = 7493.26
$ws.Cells.Item(15, 6).Value = 856.84
$ws.Cells.Item(15, 7).Value = 1250
$ws.Cells.Item(15, 8).Value = 9600.1
$ws.Cells.Item(15, 9).Value = 685.72
$ws.Cells.Item(15, 10).Value = 480.01

# Row 16
$ws.Cells.Item(16, 1).Value = 'FL6575'
$ws.Cells.Item(16, 2).Value = 'website'
$ws.Cells.Item(16, 3).Value = 'premium_economy'
$ws.Cells.Item(16, 4).Value = 6
$ws.Cells.Item(16, 5).Value = 1206.49
$ws.Cells.Item(16, 6).Value = 234.59
$ws.Cells.Item(16, 7).Value = 450
$ws.Cells.Item(16, 8).Value = 1891.08
$ws.Cells.Item(16, 9).Value = 315.18
$ws.Cells.Item(16, 10).Value = 94.55

# Row 17
$ws.Cells.Item(17, 1).Value = 'FL3756'
$ws.Cells.Item(17, 2).Value = 'app'
$ws.Cells.Item(17, 3).Value = 'economy'
$ws.Cells.Item(17, 4).Value = 38
$ws.Cells.Item(17, 5).Value = 4701.21
$ws.Cells.Item(17, 6).Value = 2036.3
$ws.Cells.Item(17, 7).Value = 2550
$ws.Cells.Item(17, 8).Value = 9287.51
$ws.Cells.Item(17, 9).Value = 244.41
$ws.Cells.Item(17, 10).Value = 464.38

# Row 18
$ws.Cells.Item(18, 1).Value = 'FL6970'
$ws.Cells.Item(18, 2).Value = 'travel_agent'
$ws.Cells.Item(18, 3).Value = 'economy'
$ws.Cells.Item(18, 4).Value = 39
$ws.Cells.Item(18, 5).Value = 35667.53
$ws.Cells.Item(18, 6).Value = 2434.26
$ws.Cells.Item(18, 7).Value = 2300
$ws.Cells.Item(18, 8).Value = 40401.79
$ws.Cells.Item(18, 9).Value = 1035.94
$ws.Cells.Item(18, 10).Value = 2020.09

# Row 19
$ws.Cells.Item(19, 1).Value = 'FL5544'
$ws.Cells.Item(19, 2).Value = 'website'
$ws.Cells.Item(19, 3).Value = 'economy'
$ws.Cells.Item(19, 4).Value = 33
$ws.Cells.Item(19, 5).Value = 8202.77
$ws.Cells.Item(19, 6).Value = 1732.45
$ws.Cells.Item(19, 7).Value = 2600
$ws.Cells.Item(19, 8).Value = 12535.22
$ws.Cells.Item(19, 9).Value = 379.86
$ws.Cells.Item(19, 10).Value = 626.76

# Row 20
$ws.Cells.Item(20, 1).Value = 'FL1573'
$ws.Cells.Item(20, 2).Value = 'website'
$ws.Cells.Item(20, 3).Value = 'business'
$ws.Cells.Item(20, 4).Value = 5
$ws.Cells.Item(20, 5).Value = 1707.25
$ws.Cells.Item(20, 6).Value = 177.01
$ws.Cells.Item(20, 7).Value = 400
$ws.Cells.Item(20, 8).Value = 2284.26
$ws.Cells.Item(20, 9).Value = 456.85
$ws.Cells.Item(20, 10).Value = 114.21

# Row 21
$ws.Cells.Item(21, 1).Value = 'FL1524'
$ws.Cells.Item(21, 2).Value = 'travel_agent'
$ws.Cells.Item(21, 3).Value = 'economy'
$ws.Cells.Item(21, 4).Value = 22
$ws.Cells.Item(21, 5).Value = 4540.28
$ws.Cells.Item(21, 6).Value = 1111.92
$ws.Cells.Item(21, 7).Value = 1300
$ws.Cells.Item(21, 8).Value = 6952.2
$ws.Cells.Item(21, 9).Value = 316.01
$ws.Cells.Item(21, 10).Value = 347.61

# Row 22
$ws.Cells.Item(22, 1).Value = 'FL7598'
$ws.Cells.Item(22, 2).Value = 'website'
$ws.Cells.Item(22, 3).Value = 'business'
$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(22, 5).Value = 8200.82
$ws.Cells.Item(22, 6).Value = 170.34
$ws.Cells.Item(22, 7).Value = 500
$ws.Cells.Item(22, 8).Value = 8871.16
$ws.Cells.Item(22, 9).Value = 1478.53
$ws.Cells.Item(22, 10).Value = 443.56

# Row 23
$ws.Cells.Item(23, 1).Value = 'FL2625'
$ws.Cells.Item(23, 2).Value = 'app'
$ws.Cells.Item(23, 3).Value = 'business'
$ws.Cells.Item(23, 4).Value = 6
$ws.Cells.Item(23, 5).Value = 1121.22
$ws.Cells.Item(23, 6).Value = 341.88
$ws.Cells.Item(23, 7).Value = 600
$ws.Cells.Item(23, 8).Value = 2063.1
$ws.Cells.Item(23, 9).Value = 343.85
$ws.Cells.Item(23, 10).Value = 103.16

# Row 24
$ws.Cells.Item(24, 1).Value = 'FL3865'
$ws.Cells.Item(24, 2).Value = 'app'
$ws.Cells.Item(24, 3).Value = 'economy'
$ws.Cells.Item(24, 4).Value = 35
$ws.Cells.Item(24, 5).Value = 8711.639999999999
$ws.Cells.Item(24, 6).Value = 1689.24
$ws.Cells.Item(24, 7).Value = 2850
$ws.Cells.Item(24, 8).Value = 13250.88
$ws.Cells.Item(24, 9).Value = 378.6
$ws.Cells.Item(24, 10).Value = 662.54

# Row 25
$ws.Cells.Item(25, 1).Value = 'FL8222'
$ws.Cells.Item(25, 2).Value = 'website'
$ws.Cells.Item(25, 3).Value = 'economy'
$ws.Cells.Item(25, 4).Value = 36
$ws.Cells.Item(25, 5).Value = 5271.08
$ws.Cells.Item(25, 6).Value = 1676.67
$ws.Cells.Item(25, 7).Value = 2150
$ws.Cells.Item(25, 8).Value = 9097.75
$ws.Cells.Item(25, 9).Value = 252.72
$ws.Cells.Item(25, 10).Value = 454.89

# Row 26
$ws.Cells.Item(26, 1).Value = 'FL6679'
$ws.Cells.Item(26, 2).Value = 'app'
$ws.Cells.Item(26, 3).Value = 'economy'
$ws.Cells.Item(26, 4).Value = 21
$ws.Cells.Item(26, 5).Value = 3035.06
$ws.Cells.Item(26, 6).Value = 999.75
$ws.Cells.Item(26, 7).Value = 1650
$ws.Cells.Item(26, 8).Value = 5684.81
$ws.Cells.Item(26, 9).Value = 270.71
$ws.Cells.Item(26, 10).Value = 284.24
